$d = $word.ActiveDocument

function Get-ParaIndex($substr) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs($i).Range.Text -like "*$substr*") {
            return $i
        }
    }
    return -1
}

# Edit 1: "... Example presentation, Sams recipe that he is making" paragraph
#   -> change "Sams" to "Sam’s" and drop the (now unnecessary) spell/grammar-check
#      proofErr markers that bracketed the misspelling.
$idx1 = Get-ParaIndex "Sams"
$p1 = $d.Paragraphs($idx1).Range
$frag1 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00DB3338" w:rsidRDefault="00D2607C"><w:pPr><w:spacing w:after="200" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Corbel" w:eastAsia="Corbel" w:hAnsi="Corbel" w:cs="Corbel"/><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Corbel" w:eastAsia="Corbel" w:hAnsi="Corbel" w:cs="Corbel"/><w:sz w:val="24"/></w:rPr><w:tab/><w:t>-</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Corbel" w:eastAsia="Corbel" w:hAnsi="Corbel" w:cs="Corbel"/><w:sz w:val="24"/></w:rPr><w:tab/><w:t xml:space="preserve">Example presentation, </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Corbel" w:eastAsia="Corbel" w:hAnsi="Corbel" w:cs="Corbel"/><w:sz w:val="24"/></w:rPr><w:t>Sam’s</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Corbel" w:eastAsia="Corbel" w:hAnsi="Corbel" w:cs="Corbel"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> recipe that he is making</w:t></w:r></w:p>'
$p1.InsertXML($frag1)

# Edit 2: "GDS and IDS reports*" -> "A software report on how iterations went"
#   (text now lives in its own run, after a run that only holds the leading tab)
$idx2 = Get-ParaIndex "GDS and IDS reports"
$p2 = $d.Paragraphs($idx2).Range
$frag2 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00DB3338" w:rsidRDefault="00D2607C"><w:pPr><w:spacing w:after="200" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Corbel" w:eastAsia="Corbel" w:hAnsi="Corbel" w:cs="Corbel"/><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Corbel" w:eastAsia="Corbel" w:hAnsi="Corbel" w:cs="Corbel"/><w:sz w:val="24"/></w:rPr><w:tab/><w:t>-</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Corbel" w:eastAsia="Corbel" w:hAnsi="Corbel" w:cs="Corbel"/><w:sz w:val="24"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii="Corbel" w:eastAsia="Corbel" w:hAnsi="Corbel" w:cs="Corbel"/><w:sz w:val="24"/></w:rPr><w:t>A software report on how iterations went</w:t></w:r></w:p>'
$p2.InsertXML($frag2)

# Edit 3: move the "_GoBack" bookmark from the end of the "Cash flow*" paragraph
#   to the end of "The Team" paragraph, two paragraphs further down (the blank
#   lastRenderedPageBreak paragraph in between is unaffected).
$idxCash = Get-ParaIndex "Cash flow"
$idxTeam = Get-ParaIndex "The Team"
$s = $d.Paragraphs($idxCash).Range.Start
$e = $d.Paragraphs($idxTeam).Range.End
$rng3 = $d.Range($s, $e)
$frag3 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00DB3338" w:rsidRDefault="00D2607C"><w:pPr><w:spacing w:after="200" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Corbel" w:eastAsia="Corbel" w:hAnsi="Corbel" w:cs="Corbel"/><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Corbel" w:eastAsia="Corbel" w:hAnsi="Corbel" w:cs="Corbel"/><w:sz w:val="24"/></w:rPr><w:tab/><w:t>-</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Corbel" w:eastAsia="Corbel" w:hAnsi="Corbel" w:cs="Corbel"/><w:sz w:val="24"/></w:rPr><w:tab/><w:t>Cash flow*</w:t></w:r></w:p><w:p w:rsidR="00DB3338" w:rsidRDefault="00D2607C"><w:pPr><w:spacing w:after="200" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Corbel" w:eastAsia="Corbel" w:hAnsi="Corbel" w:cs="Corbel"/><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Corbel" w:eastAsia="Corbel" w:hAnsi="Corbel" w:cs="Corbel"/><w:sz w:val="24"/></w:rPr><w:lastRenderedPageBreak/><w:tab/></w:r></w:p><w:p w:rsidR="00DB3338" w:rsidRDefault="00D2607C"><w:pPr><w:spacing w:after="200" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Corbel" w:eastAsia="Corbel" w:hAnsi="Corbel" w:cs="Corbel"/><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Corbel" w:eastAsia="Corbel" w:hAnsi="Corbel" w:cs="Corbel"/><w:sz w:val="24"/></w:rPr><w:t>The Team</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
$rng3.InsertXML($frag3)
